# "Add corrections and demo"
#
# 1) Remove the old "Demo: Working with files" slide (position 13).
# 2) Fix the "Storage of data" title on slide 8 so it is a single run
#    (merges the two runs "Хранение " + "данных").
# 3) The old "Demo: Using file dialogs" slide (position 23 before the
#    deletion above, position 22 afterwards) gets re-positioned and its
#    title text shortened to just "Демонстрация".

$p = $ppt.ActivePresentation

# --- 1. Delete the "Демо: Работа с файлами" slide ------------------------
$oldDemoSlide = $p.Slides.Item(13)
$oldDemoSlide.Delete()

# --- 2. "Хранение " + "данных" -> single run "Хранение данных" ----------
$storageSlide = $p.Slides.Item(8)
$storageTitle = $storageSlide.Shapes.Item(2)
$storageTitle.TextFrame.TextRange.Text = "placeholder_tmp_storage_title"
$storageTitle.TextFrame.TextRange.Text = "Хранение данных"

# --- 3. Re-layout + retitle the demo slide now at position 22 -----------
$demoSlide = $p.Slides.Item(22)
$demoTitle = $demoSlide.Shapes.Item(1)

$demoTitle.Top = 232.19692913385828
$demoTitle.Height = 75.60622047244094

$titleRange = $demoTitle.TextFrame.TextRange
# Drop "Демо" + ": " + the line break, keeping just the rest of the
# paragraph (plus its trailing endParaRPr) so the run we are left with
# carries clean formatting (no "err" typo flag) before retyping it.
$headChars = $titleRange.Characters(1, 8)
$headChars.Delete()

$titleRange.Text = "placeholder_tmp_demo_title"
$titleRange.Text = "Демонстрация"
